$d = $word.ActiveDocument

# Paragraph 1: "AI Agents in the Data Pipeline" -> Bold, size 14pt (28 half-points)
$p1 = $d.Paragraphs.Item(1).Range
$p1.Font.Bold = $true
$p1.Font.BoldBi = $true
$p1.Font.Size = 14
$p1.Font.SizeBi = 14

# Paragraph 2: "Can we build an Agentic AI Data Analyst (intern)?" -> Italic
$p2 = $d.Paragraphs.Item(2).Range
$p2.Font.Italic = $true
$p2.Font.ItalicBi = $true
